$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select B3
$ws.Range("B3").Select()

# Autofit column B to the widest content (the long navigational status text in B3)
$ws.Columns.Item(2).AutoFit() | Out-Null
